$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the "Ready for handoff" rows (4-7)
# moves from 10:37:33 to 10:37:51 (report regenerated for handoff).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-09-02 10:37:51"

# "zh-cn" sheet: Priority goes from "low" to "ht" (high) and the Latest Handoff
# Datetime is refreshed from 10:37:28 to 10:37:46 for the same four rows.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-09-02 10:37:46"

# "de-de" sheet: Priority also goes from "low" to "ht" for the same four rows,
# and its Latest Handoff Datetime shares the same refreshed timestamp as the
# Overview sheet's "Latest HO Xliff Generate Date" above.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-09-02 10:37:51"
